$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (new sentiment/price-check pass for the same ticker) ---
$ws.Cells.Item(2, 1).Value = 42651.601574074077
$ws.Cells.Item(2, 2).Value = 19
$ws.Cells.Item(2, 3).Value = "Strong Buy"
$ws.Cells.Item(2, 4).Value = 22
$ws.Cells.Item(2, 5).Value = 13716
$ws.Cells.Item(2, 6).Value = 1519
$ws.Cells.Item(2, 7).Value = 47
$ws.Cells.Item(2, 8).Value = 52
$ws.Cells.Item(2, 9).Value = 89
$ws.Cells.Item(2, 10).Value = 9
$ws.Cells.Item(2, 11).Value = 63283
$ws.Cells.Item(2, 12).Value = 153
$ws.Cells.Item(2, 13).Value = 169
$ws.Cells.Item(2, 14).Value = 74
$ws.Cells.Item(2, 15).Value = 8
$ws.Cells.Item(2, 16).Value = "Bag"
$ws.Cells.Item(2, 17).Value = 29.378539412357895
$ws.Cells.Item(2, 18).Value = 0.84
$ws.Cells.Item(2, 19).Value = -0.0136
$ws.Cells.Item(2, 20).Value = -0.03
$ws.Cells.Item(2, 21).Value = 14.53
$ws.Cells.Item(2, 22).Value = "N/A"
$ws.Cells.Item(2, 23).Value = 1
$ws.Cells.Item(2, 24).Value = 0
$ws.Cells.Item(2, 25).Value = "Up"

# --- New row 3 ---
$ws.Cells.Item(3, 1).Value = 42651.601979166669
$ws.Cells.Item(3, 1).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(3, 2).Value = 14
$ws.Cells.Item(3, 3).Value = "Buy"
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 103
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = "Bag"
$ws.Cells.Item(3, 17).Value = 29.820796582770228
$ws.Cells.Item(3, 18).Value = 0.84
$ws.Cells.Item(3, 19).Value = -0.0136
$ws.Cells.Item(3, 19).NumberFormat = "0.00%"
$ws.Cells.Item(3, 20).Value = -0.03
$ws.Cells.Item(3, 20).NumberFormat = "0.00%"
$ws.Cells.Item(3, 21).Value = 14.53
$ws.Cells.Item(3, 22).Value = "N/A"
$ws.Cells.Item(3, 23).Value = 1
$ws.Cells.Item(3, 24).Value = -0.39000000000000057
$ws.Cells.Item(3, 25).Value = "Down"

# --- New row 4 ---
$ws.Cells.Item(4, 1).Value = 42651.682523148149
$ws.Cells.Item(4, 1).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(4, 2).Value = 6
$ws.Cells.Item(4, 3).Value = "Buy"
$ws.Cells.Item(4, 4).Value = 22
$ws.Cells.Item(4, 5).Value = 13719
$ws.Cells.Item(4, 6).Value = 1519
$ws.Cells.Item(4, 7).Value = 47
$ws.Cells.Item(4, 8).Value = 52
$ws.Cells.Item(4, 9).Value = 89
$ws.Cells.Item(4, 10).Value = 9
$ws.Cells.Item(4, 11).Value = 42437
$ws.Cells.Item(4, 12).Value = 153
$ws.Cells.Item(4, 13).Value = 169
$ws.Cells.Item(4, 14).Value = 74
$ws.Cells.Item(4, 15).Value = 8
$ws.Cells.Item(4, 16).Value = "Bag"
$ws.Cells.Item(4, 17).Value = 29.009771469523784
$ws.Cells.Item(4, 18).Value = 0.84
$ws.Cells.Item(4, 19).Value = -0.0136
$ws.Cells.Item(4, 19).NumberFormat = "0.00%"
$ws.Cells.Item(4, 20).Value = -0.03
$ws.Cells.Item(4, 20).NumberFormat = "0.00%"
$ws.Cells.Item(4, 21).Value = 14.53
$ws.Cells.Item(4, 22).Value = "N/A"
$ws.Cells.Item(4, 23).Value = -2

# Column C ("Verdict") now holds the longer "Strong Buy" text, so re-fit its width
# like Excel does automatically when a cell's best-fit column receives wider content.
$ws.Columns.Item(3).AutoFit()
